$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 238, shifting the existing data (old rows
# 238-282) down to rows 242-286.
$ws.Rows("238:241").Insert()

# Populate the 4 newly inserted rows with the new weekly price records.
$newRows = @(
    @{ Row = 238; D = 45180; K = "Hass"; L = "Especial"; M = 250; N = 27000; O = 28000; P = 27600; S = 2760; T = 10 },
    @{ Row = 239; D = 45180; K = "Hass"; L = "Primera";  M = 420; N = 25000; O = 26000; P = 25476; S = 2548; T = 10 },
    @{ Row = 240; D = 45180; K = "Hass"; L = "Segunda";  M = 330; N = 23000; O = 24000; P = 23455; S = 2346; T = 10 },
    @{ Row = 241; D = 45180; K = "Hass"; L = "Tercera";  M = 290; N = 20000; O = 22000; P = 20690; S = 2069; T = 10 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 1
    $ws.Cells.Item($row, 2).Value = "Agrícola del Norte S.A. de Arica"
    $ws.Cells.Item($row, 3).Value = "Arica y Parinacota"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = 15
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100106
    $ws.Cells.Item($row, 8).Value = "Oleaginosos"
    $ws.Cells.Item($row, 9).Value = 100106002
    $ws.Cells.Item($row, 10).Value = "Palta"
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = "$/bandeja 10 kilos"
    $ws.Cells.Item($row, 18).Value = "Perú"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
